$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.665.88'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.328.89'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '237.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.660'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.45%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '71.53'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -5.88%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -7.41%  '
$ws.Range('E10').Value = '  -4.88%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.33'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '32.19'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.107'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('E14').Value = '  -6.06%  '
$ws.Range('D15').Value = '2.675.60'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.99'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -5.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.888'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.37%  '
$ws.Range('D18').Value = '2.332.27'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = '43.595.60'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('D20').Value = '0.0₃0999'
$ws.Range('E20').Value = '  -3.36%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '77.62'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.59'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '249.42'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.98%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.89'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.05%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.71'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.47'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.59%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.24'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -9.18%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.25'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '174.32'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '21.94'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.49%  '
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.133'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.32'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.02'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -5.35%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.71'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('E38').Value = '  -3.73%  '
$ws.Range('E39').Value = '  -2.90%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.62'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +25.33%  '
$ws.Range('E41').Value = '  -3.68%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '64.60'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +18.72%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '9.13'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.104'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.08%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '18.62'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.52%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.193'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -4.09%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.42'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.27%  '
$ws.Range('B49').Value = 'TrustWalletToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.21'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.21%  '
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '97.12'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.93%  '
